# Update "江西-漫展信息.xlsx" workbook output (gh-pages generated data refresh)
#
# Applies the same set of edits to both the "展览" sheet (index 1) and the
# "全部类型" sheet (index 4), since both contain the same event rows (the
# latter additionally contains a local-life/"Kpop" row inserted earlier,
# which shifts all of its row numbers by +1 relative to "展览").

$wb = $excel.ActiveWorkbook

function Update-EventSheet($ws, $offset) {
    # $offset is 0 for "展览", 1 for "全部类型"

    # --- simple value updates (row numbers below are for "展览"; add $offset) ---

    # Row 2: only-event got cancelled, lowest price becomes "不可售" (not for sale)
    $ws.Range("C" + (2 + $offset)).Value = "南昌·原X穹X崩only（取消）"
    $ws.Range("G" + (2 + $offset)).Value = "不可售"

    # "想去人数" (interest count) refreshes for a handful of still-open events
    $ws.Range("F" + (3 + $offset)).Value = 162
    $ws.Range("F" + (4 + $offset)).Value = 167
    $ws.Range("F" + (5 + $offset)).Value = 4774
    $ws.Range("F" + (8 + $offset)).Value = 526
    $ws.Range("F" + (12 + $offset)).Value = 1348
    $ws.Range("F" + (13 + $offset)).Value = 3022
    $ws.Range("F" + (14 + $offset)).Value = 389
    $ws.Range("F" + (15 + $offset)).Value = 100
    $ws.Range("F" + (16 + $offset)).Value = 88
    $ws.Range("F" + (18 + $offset)).Value = 2419
    $ws.Range("F" + (22 + $offset)).Value = 169

    # --- insert a brand-new event row right after the current row 22 ---
    $newRow = 23 + $offset
    $ws.Rows.Item($newRow).Insert()

    # match the bold/centered/bordered formatting used by the rest of column A
    $aboveCell = $ws.Range("A" + ($newRow - 1))
    $newCell = $ws.Range("A" + $newRow)
    $newCell.Font.Bold = $aboveCell.Font.Bold()
    $newCell.HorizontalAlignment = $aboveCell.HorizontalAlignment()
    $newCell.VerticalAlignment = $aboveCell.VerticalAlignment()
    $newCell.Borders.LineStyle = $aboveCell.Borders.LineStyle()

    $newCell.Value = 22
    $ws.Range("B" + $newRow).Value = "2024-05-03"
    $ws.Range("C" + $newRow).Value = "赣州·漫库书店次元漫展"
    $ws.Range("D" + $newRow).Value = "南门口地一大道下沉广场 漫库书店"
    $ws.Range("E" + $newRow).Value = "2024.05.03 10:00-05.04 18:00"
    $ws.Range("F" + $newRow).Value = 1
    $ws.Range("G" + $newRow).Value = 40
    $ws.Range("H" + $newRow).Value = "https://show.bilibili.com/platform/detail.html?id=83855"
    $ws.Range("I" + $newRow).Value = "//i0.hdslb.com/bfs/openplatform/202404/juDVRy6Y1712481590113.jpeg"

    # --- the old "南昌·代号鸢盛花行only" row (now pushed down one row by the
    #     insert above) also gets a refreshed interest count ---
    $ws.Range("F" + (26 + $offset)).Value = 242
}

# Sheet 1: "展览"
$wsExhibition = $wb.Worksheets.Item(1)
Update-EventSheet $wsExhibition 0

# Sheet 4: "全部类型"
$wsAll = $wb.Worksheets.Item(4)
Update-EventSheet $wsAll 1
